$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.169212666666667
$ws.Range("H2").Value = 3.507638
$ws.Range("I2").Value = 0.0005193657195729173
$ws.Range("J2").Value = 0.0005193657195729173
$ws.Range("M2").Value = 0.003643333333333333
$ws.Range("N2").Value = 0.01093
$ws.Range("O2").Value = 0.002177035403614994
$ws.Range("P2").Value = 0.002177035403614994
$ws.Range("Q2").Value = 0.004259831482222223
$ws.Range("R2").Value = 0.03833848334
$ws.Range("S2").Value = 0.000001130677558934218
$ws.Range("T2").Value = 0.000001130677558934218

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.169212666666667
$ws.Range("H3").Value = 3.507638
$ws.Range("I3").Value = 0.0005193657195729173
$ws.Range("J3").Value = 0.0005193657195729173
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.669886333333333
$ws.Range("N3").Value = 5.009659
$ws.Range("O3").Value = 0.997822964596385
$ws.Range("P3").Value = 0.997822964596385
$ws.Range("Q3").Value = 1.952452252826889
$ws.Range("R3").Value = 17.572070275442
$ws.Range("S3").Value = 0.0005182350420139831
$ws.Range("T3").Value = 0.0005182350420139831

$ws.Range("I4").Value = 0.9638330474556795
$ws.Range("J4").Value = 0.9638330474556795
$ws.Range("M4").Value = 0.003643333333333333
$ws.Range("N4").Value = 0.01093
$ws.Range("O4").Value = 0.002177035403614994
$ws.Range("P4").Value = 0.002177035403614994
$ws.Range("Q4").Value = 7.90534724266
$ws.Range("R4").Value = 71.14812518394001
$ws.Range("S4").Value = 0.002098298667485145
$ws.Range("T4").Value = 0.002098298667485145

$ws.Range("I5").Value = 0.9638330474556795
$ws.Range("J5").Value = 0.9638330474556795
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.669886333333333
$ws.Range("N5").Value = 5.009659
$ws.Range("O5").Value = 0.997822964596385
$ws.Range("P5").Value = 0.997822964596385
$ws.Range("Q5").Value = 3623.338880358358
$ws.Range("R5").Value = 32610.04992322522
$ws.Range("S5").Value = 0.9617347487881944
$ws.Range("T5").Value = 0.9617347487881944

$ws.Range("G6").Value = 80.250984
$ws.Range("H6").Value = 240.752952
$ws.Range("I6").Value = 0.03564758682474761
$ws.Range("J6").Value = 0.0356475868247476
$ws.Range("M6").Value = 0.003643333333333333
$ws.Range("N6").Value = 0.01093
$ws.Range("O6").Value = 0.002177035403614994
$ws.Range("P6").Value = 0.002177035403614994
$ws.Range("Q6").Value = 0.29238108504
$ws.Range("R6").Value = 2.63142976536
$ws.Range("S6").Value = 0.00007760605857091495
$ws.Range("T6").Value = 0.00007760605857091494

$ws.Range("G7").Value = 80.250984
$ws.Range("H7").Value = 240.752952
$ws.Range("I7").Value = 0.03564758682474761
$ws.Range("J7").Value = 0.0356475868247476
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.669886333333333
$ws.Range("N7").Value = 5.009659
$ws.Range("O7").Value = 0.997822964596385
$ws.Range("P7").Value = 0.997822964596385
$ws.Range("Q7").Value = 134.010021418152
$ws.Range("R7").Value = 1206.090192763368
$ws.Range("S7").Value = 0.03556998076617669
$ws.Range("T7").Value = 0.03556998076617669

$wb.Save()
